# Applies the RAD test-case / MRF data update described in the commit:
# "Added RAD Test Cases and data for MRF."
#
# For every worksheet this rewrites the Result (A) / Date (B) columns of the
# test-run rows with a fresh batch of Pass/Fail results and run timestamps,
# adding rows that previously had no Result/Date recorded, and moves the
# active selection / active tab to where Excel left it after the edit.

$wb = $excel.ActiveWorkbook

function Set-ResultDate {
    param(
        $ws,
        [int]$row,
        [string]$result,
        [string]$date
    )
    if ($result -ne $null) {
        $cellA = $ws.Cells.Item($row, 1)
        $cellA.Value = $result
        $cellA.Style = "Normal"
    }
    $cellB = $ws.Cells.Item($row, 2)
    $cellB.Value = $date
    $cellB.Style = "Normal"
}

# ---------------------------------------------------------------------------
# Sheet "Estimated"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Estimated")
Set-ResultDate $ws 2 "Pass" "Mon Oct 02 18:35:48 EDT 2023"
Set-ResultDate $ws 3 "Pass" "Mon Oct 02 18:36:37 EDT 2023"
Set-ResultDate $ws 4 "Fail" "Mon Oct 02 18:37:23 EDT 2023"
Set-ResultDate $ws 5 "Fail" "Mon Oct 02 18:38:21 EDT 2023"
Set-ResultDate $ws 6 "Fail" "Mon Oct 02 18:39:21 EDT 2023"
Set-ResultDate $ws 7 "Fail" "Mon Oct 02 18:40:19 EDT 2023"
$ws.Range("D4:E7").Select()

# ---------------------------------------------------------------------------
# Sheet "Existing"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Existing")
Set-ResultDate $ws 2 $null "Mon Oct 02 18:43:15 EDT 2023"
Set-ResultDate $ws 3 $null "Mon Oct 02 18:44:00 EDT 2023"
Set-ResultDate $ws 4 $null "Mon Oct 02 18:44:43 EDT 2023"
Set-ResultDate $ws 5 $null "Mon Oct 02 18:45:27 EDT 2023"
Set-ResultDate $ws 6 $null "Mon Oct 02 18:46:10 EDT 2023"
Set-ResultDate $ws 7 $null "Mon Oct 02 18:46:53 EDT 2023"
Set-ResultDate $ws 8 $null "Mon Oct 02 18:47:37 EDT 2023"
Set-ResultDate $ws 9 $null "Mon Oct 02 18:48:20 EDT 2023"
Set-ResultDate $ws 10 $null "Mon Oct 02 18:49:04 EDT 2023"
Set-ResultDate $ws 11 $null "Mon Oct 02 18:49:47 EDT 2023"
Set-ResultDate $ws 12 $null "Mon Oct 02 18:50:45 EDT 2023"

# ---------------------------------------------------------------------------
# Sheet "Extension"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Extension")
Set-ResultDate $ws 2 "Pass" "Mon Oct 02 18:51:45 EDT 2023"
Set-ResultDate $ws 3 "Pass" "Mon Oct 02 18:52:29 EDT 2023"
Set-ResultDate $ws 4 $null "Mon Oct 02 18:53:12 EDT 2023"
Set-ResultDate $ws 5 $null "Mon Oct 02 18:54:10 EDT 2023"
Set-ResultDate $ws 6 $null "Mon Oct 02 18:55:07 EDT 2023"
Set-ResultDate $ws 7 $null "Mon Oct 02 18:56:05 EDT 2023"
$ws.Range("D4:E7").Select()

# ---------------------------------------------------------------------------
# Sheet "NewTaxReturn"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("NewTaxReturn")
Set-ResultDate $ws 2 "Pass" "Mon Oct 02 18:57:04 EDT 2023"
Set-ResultDate $ws 3 $null "Mon Oct 02 18:57:47 EDT 2023"
Set-ResultDate $ws 4 $null "Mon Oct 02 18:58:45 EDT 2023"
Set-ResultDate $ws 5 $null "Mon Oct 02 18:59:43 EDT 2023"
Set-ResultDate $ws 6 $null "Mon Oct 02 19:00:42 EDT 2023"
Set-ResultDate $ws 7 "Pass" "Mon Oct 02 19:01:41 EDT 2023"
Set-ResultDate $ws 8 $null "Mon Oct 02 19:02:23 EDT 2023"
Set-ResultDate $ws 9 $null "Mon Oct 02 19:03:21 EDT 2023"
Set-ResultDate $ws 10 $null "Mon Oct 02 19:04:20 EDT 2023"
Set-ResultDate $ws 11 $null "Mon Oct 02 19:05:18 EDT 2023"
Set-ResultDate $ws 12 "Pass" "Mon Oct 02 19:06:16 EDT 2023"
Set-ResultDate $ws 13 $null "Mon Oct 02 19:06:59 EDT 2023"
Set-ResultDate $ws 14 $null "Mon Oct 02 19:07:58 EDT 2023"
Set-ResultDate $ws 15 $null "Mon Oct 02 19:08:56 EDT 2023"
Set-ResultDate $ws 16 $null "Mon Oct 02 19:09:54 EDT 2023"
$ws.Range("B20").Select()
$ws.Activate()

# ---------------------------------------------------------------------------
# Sheet "Personal_EL"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Personal_EL")
Set-ResultDate $ws 2 $null "Mon Oct 02 19:10:53 EDT 2023"

# ---------------------------------------------------------------------------
# Sheet "Personal_IND"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Personal_IND")
Set-ResultDate $ws 2 $null "Mon Oct 02 19:11:35 EDT 2023"
Set-ResultDate $ws 3 $null "Mon Oct 02 19:12:15 EDT 2023"
Set-ResultDate $ws 4 "Pass" "Mon Oct 02 19:12:55 EDT 2023"
Set-ResultDate $ws 5 "Pass" "Mon Oct 02 19:13:35 EDT 2023"
Set-ResultDate $ws 6 "Pass" "Mon Oct 02 19:14:15 EDT 2023"

# ---------------------------------------------------------------------------
# Sheet "Personal_JNT"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Personal_JNT")
Set-ResultDate $ws 2 $null "Mon Oct 02 19:14:56 EDT 2023"
Set-ResultDate $ws 3 $null "Mon Oct 02 19:15:43 EDT 2023"
Set-ResultDate $ws 4 "Pass" "Mon Oct 02 19:16:30 EDT 2023"
Set-ResultDate $ws 5 "Pass" "Mon Oct 02 19:17:17 EDT 2023"
Set-ResultDate $ws 6 "Pass" "Mon Oct 02 19:18:03 EDT 2023"
